{"js": "// Remove the empty paragraph (leftover \"_GoBack\" bookmark paragraph with no\n// text) that sits between \"A simple demonstration of a query :\" and\n// \"End of demonstration.\" \u2014 this is the blank line an empty AQL expression\n// used to generate (see #418).\n\nconst bookmarkRange = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nbookmarkRange.load(\"isNullObject\");\nawait context.sync();\n\nif (!bookmarkRange.isNullObject) {\n  const paragraphs = bookmarkRange.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  paragraphs.items[0].delete();\n  await context.sync();\n}\n", "ps1": "# Remove the empty paragraph (leftover \"_GoBack\" bookmark paragraph with no\n# text) that sits between \"A simple demonstration of a query :\" and\n# \"End of demonstration.\" \u2014 this is the blank line an empty AQL expression\n# used to generate (see #418).\n\n$d = $word.ActiveDocument\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $bm = $d.Bookmarks.Item(\"_GoBack\")\n    $bmStart = $bm.Start\n\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($bmStart -ge $p.Range.Start -and $bmStart -lt $p.Range.End) {\n            $p.Range.Delete()\n            break\n        }\n    }\n}\n"}
